# Auto-generated script to update TPM values in Gnai2-Cxcr1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 181.4944075
$ws.Range("H2").Value = 362.988815
$ws.Range("I2").Value = 0.2239486468210351
$ws.Range("J2").Value = 0.1654349085470023
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3657716666666667
$ws.Range("N2").Value = 1.097315
$ws.Range("O2").Value = 0.03059585711603819
$ws.Range("P2").Value = 0.03059585711603819
$ws.Range("Q2").Value = 66.38551192195418
$ws.Range("R2").Value = 398.3130715317251
$ws.Range("S2").Value = 0.00685190079946649
$ws.Range("T2").Value = 0.005061622823908928
# Row 3
$ws.Range("G3").Value = 181.4944075
$ws.Range("H3").Value = 362.988815
$ws.Range("I3").Value = 0.2239486468210351
$ws.Range("J3").Value = 0.1654349085470023
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.510206
$ws.Range("N3").Value = 34.530618
$ws.Range("O3").Value = 0.9627990635838353
$ws.Range("P3").Value = 0.9627990635838353
$ws.Range("Q3").Value = 2089.038018172945
$ws.Range("R3").Value = 12534.22810903767
$ws.Range("S3").Value = 0.2156175474501597
$ws.Range("T3").Value = 0.1592805750331313
# Row 4
$ws.Range("G4").Value = 181.4944075
$ws.Range("H4").Value = 362.988815
$ws.Range("I4").Value = 0.2239486468210351
$ws.Range("J4").Value = 0.1654349085470023
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.07896333333333333
$ws.Range("N4").Value = 0.23689
$ws.Range("O4").Value = 0.006605079300126477
$ws.Range("P4").Value = 0.006605079300126477
$ws.Range("Q4").Value = 14.33140339755833
$ws.Range("R4").Value = 85.98842038535001
$ws.Range("S4").Value = 0.001479198571408954
$ws.Range("T4").Value = 0.001092710689962122
# Row 5
$ws.Range("I5").Value = 0.07700606288633029
$ws.Range("J5").Value = 0.08532865336765341
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3657716666666667
$ws.Range("N5").Value = 1.097315
$ws.Range("O5").Value = 0.03059585711603819
$ws.Range("P5").Value = 0.03059585711603819
$ws.Range("Q5").Value = 22.82704976506722
$ws.Range("R5").Value = 205.443447885605
$ws.Range("S5").Value = 0.002356066497138813
$ws.Range("T5").Value = 0.002610703286340675
# Row 6
$ws.Range("I6").Value = 0.07700606288633029
$ws.Range("J6").Value = 0.08532865336765341
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.510206
$ws.Range("N6").Value = 34.530618
$ws.Range("O6").Value = 0.9627990635838353
$ws.Range("P6").Value = 0.9627990635838353
$ws.Range("Q6").Value = 718.3280420886673
$ws.Range("R6").Value = 6464.952378798006
$ws.Range("S6").Value = 0.07414136523723673
$ws.Range("T6").Value = 0.08215434755924637
# Row 7
$ws.Range("I7").Value = 0.07700606288633029
$ws.Range("J7").Value = 0.08532865336765341
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.07896333333333333
$ws.Range("N7").Value = 0.23689
$ws.Range("O7").Value = 0.006605079300126477
$ws.Range("P7").Value = 0.006605079300126477
$ws.Range("Q7").Value = 4.927937573847778
$ws.Range("R7").Value = 44.35143816463
$ws.Range("S7").Value = 0.0005086311519547379
$ws.Range("T7").Value = 0.000563602522066355
# Row 8
$ws.Range("G8").Value = 171.9980316666667
$ws.Range("H8").Value = 515.994095
$ws.Range("I8").Value = 0.2122309275432167
$ws.Range("J8").Value = 0.235168226649403
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3657716666666667
$ws.Range("N8").Value = 1.097315
$ws.Range("O8").Value = 0.03059585711603819
$ws.Range("P8").Value = 0.03059585711603819
$ws.Range("Q8").Value = 62.91200670610277
$ws.Range("R8").Value = 566.208060354925
$ws.Range("S8").Value = 0.00649338713471651
$ws.Range("T8").Value = 0.007195173460797218
# Row 9
$ws.Range("G9").Value = 171.9980316666667
$ws.Range("H9").Value = 515.994095
$ws.Range("I9").Value = 0.2122309275432167
$ws.Range("J9").Value = 0.235168226649403
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.510206
$ws.Range("N9").Value = 34.530618
$ws.Range("O9").Value = 0.9627990635838353
$ws.Range("P9").Value = 0.9627990635838353
$ws.Range("Q9").Value = 1979.732776077856
$ws.Range("R9").Value = 17817.59498470071
$ws.Range("S9").Value = 0.2043357383021378
$ws.Range("T9").Value = 0.2264197484027164
# Row 10
$ws.Range("G10").Value = 171.9980316666667
$ws.Range("H10").Value = 515.994095
$ws.Range("I10").Value = 0.2122309275432167
$ws.Range("J10").Value = 0.235168226649403
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.07896333333333333
$ws.Range("N10").Value = 0.23689
$ws.Range("O10").Value = 0.006605079300126477
$ws.Range("P10").Value = 0.006605079300126477
$ws.Range("Q10").Value = 13.58153790717222
$ws.Range("R10").Value = 122.23384116455
$ws.Range("S10").Value = 0.001401802106362342
$ws.Range("T10").Value = 0.001553304785889424
# Row 11
$ws.Range("G11").Value = 55.64279550000001
$ws.Range("H11").Value = 111.285591
$ws.Range("I11").Value = 0.06865847234198982
$ws.Range("J11").Value = 0.05071925307032974
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3657716666666667
$ws.Range("N11").Value = 1.097315
$ws.Range("O11").Value = 0.03059585711603819
$ws.Range("P11").Value = 0.03059585711603819
$ws.Range("Q11").Value = 20.3525580480275
$ws.Range("R11").Value = 122.115348288165
$ws.Range("S11").Value = 0.00210066480958098
$ws.Range("T11").Value = 0.00155179901997199
# Row 12
$ws.Range("G12").Value = 55.64279550000001
$ws.Range("H12").Value = 111.285591
$ws.Range("I12").Value = 0.06865847234198982
$ws.Range("J12").Value = 0.05071925307032974
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 11.510206
$ws.Range("N12").Value = 34.530618
$ws.Range("O12").Value = 0.9627990635838353
$ws.Range("P12").Value = 0.9627990635838353
$ws.Range("Q12").Value = 640.460038620873
$ws.Range("R12").Value = 3842.760231725238
$ws.Range("S12").Value = 0.06610431287796445
$ws.Range("T12").Value = 0.04883244936178504
# Row 13
$ws.Range("G13").Value = 55.64279550000001
$ws.Range("H13").Value = 111.285591
$ws.Range("I13").Value = 0.06865847234198982
$ws.Range("J13").Value = 0.05071925307032974
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.07896333333333333
$ws.Range("N13").Value = 0.23689
$ws.Range("O13").Value = 0.006605079300126477
$ws.Range("P13").Value = 0.006605079300126477
$ws.Range("Q13").Value = 4.393740608665
$ws.Range("R13").Value = 26.36244365199
$ws.Range("S13").Value = 0.0004534946544443833
$ws.Range("T13").Value = 0.0003350046885727112
# Row 14
$ws.Range("G14").Value = 203.386317
$ws.Range("H14").Value = 610.158951
$ws.Range("I14").Value = 0.250961399315095
$ws.Range("J14").Value = 0.2780845747487284
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.3657716666666667
$ws.Range("N14").Value = 1.097315
$ws.Range("O14").Value = 0.03059585711603819
$ws.Range("P14").Value = 0.03059585711603819
$ws.Range("Q14").Value = 74.392952146285
$ws.Range("R14").Value = 669.5365693165651
$ws.Range("S14").Value = 0.007678379115085651
$ws.Range("T14").Value = 0.008508235915186336
# Row 15
$ws.Range("G15").Value = 203.386317
$ws.Range("H15").Value = 610.158951
$ws.Range("I15").Value = 0.250961399315095
$ws.Range("J15").Value = 0.2780845747487284
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 11.510206
$ws.Range("N15").Value = 34.530618
$ws.Range("O15").Value = 0.9627990635838353
$ws.Range("P15").Value = 0.9627990635838353
$ws.Range("Q15").Value = 2341.018406251302
$ws.Range("R15").Value = 21069.16565626172
$ws.Range("S15").Value = 0.2416254002562624
$ws.Range("T15").Value = 0.2677395681651848
# Row 16
$ws.Range("G16").Value = 203.386317
$ws.Range("H16").Value = 610.158951
$ws.Range("I16").Value = 0.250961399315095
$ws.Range("J16").Value = 0.2780845747487284
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.07896333333333333
$ws.Range("N16").Value = 0.23689
$ws.Range("O16").Value = 0.006605079300126477
$ws.Range("P16").Value = 0.006605079300126477
$ws.Range("Q16").Value = 16.06006154471
$ws.Range("R16").Value = 144.54055390239
$ws.Range("S16").Value = 0.001657619943746909
$ws.Range("T16").Value = 0.0018367706683573
# Row 17
$ws.Range("G17").Value = 135.4992116666667
$ws.Range("H17").Value = 406.497635
$ws.Range("I17").Value = 0.167194491092333
$ws.Range("J17").Value = 0.1852643836168829
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3657716666666667
$ws.Range("N17").Value = 1.097315
$ws.Range("O17").Value = 0.03059585711603819
$ws.Range("P17").Value = 0.03059585711603819
$ws.Range("Q17").Value = 49.56177248333611
$ws.Range("R17").Value = 446.055952350025
$ws.Range("S17").Value = 0.00511545876004974
$ws.Range("T17").Value = 0.005668322609833034
# Row 18
$ws.Range("G18").Value = 135.4992116666667
$ws.Range("H18").Value = 406.497635
$ws.Range("I18").Value = 0.167194491092333
$ws.Range("J18").Value = 0.1852643836168829
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 11.510206
$ws.Range("N18").Value = 34.530618
$ws.Range("O18").Value = 0.9627990635838353
$ws.Range("P18").Value = 0.9627990635838353
$ws.Range("Q18").Value = 1559.623839120937
$ws.Range("R18").Value = 14036.61455208843
$ws.Range("S18").Value = 0.1609746994600741
$ws.Range("T18").Value = 0.1783723750617713
# Row 19
$ws.Range("G19").Value = 135.4992116666667
$ws.Range("H19").Value = 406.497635
$ws.Range("I19").Value = 0.167194491092333
$ws.Range("J19").Value = 0.1852643836168829
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.07896333333333333
$ws.Range("N19").Value = 0.23689
$ws.Range("O19").Value = 0.006605079300126477
$ws.Range("P19").Value = 0.006605079300126477
$ws.Range("Q19").Value = 10.69946941723889
$ws.Range("R19").Value = 96.29522475515
$ws.Range("S19").Value = 0.001104332872209149
$ws.Range("T19").Value = 0.001223685945278564

Write-Host "Done updating cells"